$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "stg_examples_images" right after the
#    "stg_examples" sheet (and before "stg_uri_pages").
# ------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("stg_examples")
$newSheet = $wb.Worksheets.Add($null, $srcSheet)
$newSheet.Name = "stg_examples_images"

# ------------------------------------------------------------------
# 2. Fill in the header row + data rows describing the new
#    stg_examples_images staging table schema.
# ------------------------------------------------------------------
$headers = @("table_name", "field_name", "data_type", "primary_key", "source")
for ($col = 1; $col -le $headers.Length; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$rows = @(
    @("stg_examples_images", "examples_images_id", "SERIAL", $true, "Database generated"),
    @("stg_examples_images", "stg_examples_id", "INTEGER", $false, "table: stg_examples, field: examples_id"),
    @("stg_examples_images", "image_type", "VARCHAR", $false, "One of: obverse, reverse, both sides, unknown"),
    @("stg_examples_images", "link", "VARCHAR", $false, "Scraping URI page"),
    @("stg_examples_images", "ts", "TIMESTAMP", $false, "Database generated")
)

$r = 2
foreach ($row in $rows) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $newSheet.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}

# ------------------------------------------------------------------
# 3. Turn the range into an Excel Table (ListObject) like the other
#    "schema description" sheets in this workbook.
# ------------------------------------------------------------------
$lastRow = 1 + $rows.Length
$tableRange = $newSheet.Range("A1:E" + $lastRow)
$lo = $newSheet.ListObjects.Add(1, $tableRange, $null, 1)

# ------------------------------------------------------------------
# 4. Match the column widths used by the new sheet.
# ------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 17.1640625
$newSheet.Columns.Item(2).ColumnWidth = 22.83203125
$newSheet.Columns.Item(3).ColumnWidth = 17.33203125
$newSheet.Columns.Item(4).ColumnWidth = 17.33203125
$newSheet.Columns.Item(5).ColumnWidth = 21.83203125

# ------------------------------------------------------------------
# 5. Update the "stg_examples" sheet: the source note for the
#    uri_examples_id field now clarifies it comes from the URI page.
# ------------------------------------------------------------------
$examplesSheet = $wb.Worksheets.Item("stg_examples")
$examplesSheet.Range("E11").Value = "Scraping URI page. (This is from URI page)"

# ------------------------------------------------------------------
# 6. Leave the new sheet selected/active, matching the captured view.
# ------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("B4").Select()
